$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header row becomes generic column headers
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Acceptance Criteria"
$ws.Range("D1").Value = "Priority"
$ws.Range("E1").Value = "Story Points"
$ws.Range("F1").Value = "Dependencies"
$ws.Range("G1").Value = "Notes"

# Row 2 - Unified Operations Dashboard story
$ws.Range("A2").Value = "Design Unified Operations Dashboard Wireframes and Mockups"
$ws.Range("B2").Value = "As a manager, I want to view a unified operations dashboard so that I can access all core modules clearly and intuitively."
$ws.Range("C2").Value = "Feature: Unified Operations Dashboard`n  Scenario: Manager views dashboard`n    Given a manager logs in`n    When they view their dashboard`n    Then all core modules are presented clearly and intuitively`n    And the design is consistent with the established brand and style guide"
$ws.Range("D2").Value = "High"
$ws.Range("E2").Value = "'8"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "Wireframes and mockups should cover all dashboard modules and follow the brand style guide."

# Row 3 - Incident Reporting Form story
$ws.Range("A3").Value = "Design Integrated Incident Reporting Form and Submission Flow"
$ws.Range("B3").Value = "As a staff member, I want to report an incident using a simple and quick form so that I can complete the process with minimal training."
$ws.Range("C3").Value = "Feature: Incident Reporting Form`n  Scenario: Staff member reports an incident`n    Given a staff member needs to report an incident`n    When they access the form`n    Then the process is simple, quick, and requires minimal training`n    And the design is consistent with the established brand and style guide"
$ws.Range("D3").Value = "High"
$ws.Range("E3").Value = "'5"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "Form design should minimize required fields and steps, and be intuitive for first-time users."

# Row 4 - Staff Task & Schedule Viewer story
$ws.Range("A4").Value = "Create UI for Staff Task & Schedule Viewer"
$ws.Range("B4").Value = "As a staff member, I want to view my tasks and schedule in a dedicated UI so that I can easily manage my responsibilities."
$ws.Range("C4").Value = "Feature: Staff Task & Schedule Viewer`n  Scenario: Staff member views tasks and schedule`n    Given a user views the task & schedule screen`n    When they interact with it`n    Then the design is consistent with the established brand and style guide"
$ws.Range("D4").Value = "Medium"
$ws.Range("E4").Value = "'5"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "UI should display tasks and schedules in a clear, organized manner, supporting filtering and sorting."

# Row 5 - Compliance & Certification Tracking Module story
$ws.Range("A5").Value = "Visual Design for Compliance & Certification Tracking Module"
$ws.Range("B5").Value = "As a user, I want to track compliance and certifications in a visually clear module so that I can monitor requirements and expirations easily."
$ws.Range("C5").Value = "Feature: Compliance & Certification Tracking`n  Scenario: User views compliance tracking module`n    Given a user views the compliance & certification tracking screen`n    When they interact with it`n    Then the design is consistent with the established brand and style guide"
$ws.Range("D5").Value = "Medium"
$ws.Range("E5").Value = "'5"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "Visual design should highlight upcoming expirations and compliance status, using brand colors and styles."

# Row 6 - new row: Closed-Loop Incident-to-Training Workflow story
$ws.Range("A6").Value = "Map User Flow for Closed-Loop Incident-to-Training Workflow"
$ws.Range("B6").Value = "As a product team member, I want a mapped user flow for the closed-loop incident-to-training workflow so that users can seamlessly transition from incident reporting to training assignment."
$ws.Range("C6").Value = "Feature: Closed-Loop Incident-to-Training Workflow`n  Scenario: User follows incident-to-training workflow`n    Given a user completes an incident report`n    When they are assigned follow-up training`n    Then the workflow is mapped clearly and is consistent with the established brand and style guide"
$ws.Range("D6").Value = "Medium"
$ws.Range("E6").Value = "'8"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "User flow mapping should include all steps from incident submission to training completion, with visual cues and guidance."
